$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ptprc"
$ws.Range("C2").Value = "Mrc1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 243.7171333333333
$ws.Range("H2").Value = 731.1514
$ws.Range("I2").Value = 0.9993032963424349
$ws.Range("J2").Value = 0.999303296342435
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.595354
$ws.Range("N2").Value = 37.786062
$ws.Range("O2").Value = 0.989145645632832
$ws.Range("P2").Value = 0.9891456456328321
$ws.Range("Q2").Value = 3069.703570198533
$ws.Range("R2").Value = 27627.3321317868
$ws.Range("S2").Value = 0.988456504243655
$ws.Range("T2").Value = 0.9884565042436553

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ptprc"
$ws.Range("C3").Value = "Mrc1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 243.7171333333333
$ws.Range("H3").Value = 731.1514
$ws.Range("I3").Value = 0.9993032963424349
$ws.Range("J3").Value = 0.999303296342435
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.07889166666666667
$ws.Range("N3").Value = 0.236675
$ws.Range("O3").Value = 0.006195566123830276
$ws.Range("P3").Value = 0.006195566123830277
$ws.Range("Q3").Value = 19.22725084388889
$ws.Range("R3").Value = 173.045257595
$ws.Range("S3").Value = 0.006191249650251117
$ws.Range("T3").Value = 0.006191249650251119

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ptprc"
$ws.Range("C4").Value = "Mrc1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 243.7171333333333
$ws.Range("H4").Value = 731.1514
$ws.Range("I4").Value = 0.9993032963424349
$ws.Range("J4").Value = 0.999303296342435
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05932299999999999
$ws.Range("N4").Value = 0.177969
$ws.Range("O4").Value = 0.004658788243337701
$ws.Range("P4").Value = 0.004658788243337702
$ws.Range("Q4").Value = 14.45803150073333
$ws.Range("R4").Value = 130.1222835066
$ws.Range("S4").Value = 0.004655542448528747
$ws.Range("T4").Value = 0.004655542448528747

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ptprc"
$ws.Range("C5").Value = "Mrc1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1280236666666667
$ws.Range("H5").Value = 0.384071
$ws.Range("I5").Value = 0.0005249301530839377
$ws.Range("J5").Value = 0.0005249301530839377
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.595354
$ws.Range("N5").Value = 37.786062
$ws.Range("O5").Value = 0.989145645632832
$ws.Range("P5").Value = 0.9891456456328321
$ws.Range("Q5").Value = 1.612503402044667
$ws.Range("R5").Value = 14.512530618402
$ws.Range("S5").Value = 0.0005192323751843529
$ws.Range("T5").Value = 0.000519232375184353

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ptprc"
$ws.Range("C6").Value = "Mrc1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1280236666666667
$ws.Range("H6").Value = 0.384071
$ws.Range("I6").Value = 0.0005249301530839377
$ws.Range("J6").Value = 0.0005249301530839377
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.07889166666666667
$ws.Range("N6").Value = 0.236675
$ws.Range("O6").Value = 0.006195566123830276
$ws.Range("P6").Value = 0.006195566123830277
$ws.Range("Q6").Value = 0.01010000043611111
$ws.Range("R6").Value = 0.090900003925
$ws.Range("S6").Value = 0.000003252239473823886
$ws.Range("T6").Value = 0.000003252239473823886

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ptprc"
$ws.Range("C7").Value = "Mrc1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1280236666666667
$ws.Range("H7").Value = 0.384071
$ws.Range("I7").Value = 0.0005249301530839377
$ws.Range("J7").Value = 0.0005249301530839377
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.05932299999999999
$ws.Range("N7").Value = 0.177969
$ws.Range("O7").Value = 0.004658788243337701
$ws.Range("P7").Value = 0.004658788243337702
$ws.Range("Q7").Value = 0.007594747977666666
$ws.Range("R7").Value = 0.06835273179899999
$ws.Range("S7").Value = 0.000002445538425760909
$ws.Range("T7").Value = 0.000002445538425760909

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ptprc"
$ws.Range("C8").Value = "Mrc1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.04189333333333333
$ws.Range("H8").Value = 0.12568
$ws.Range("I8").Value = 0.0001717735044811748
$ws.Range("J8").Value = 0.0001717735044811748
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 12.595354
$ws.Range("N8").Value = 37.786062
$ws.Range("O8").Value = 0.989145645632832
$ws.Range("P8").Value = 0.9891456456328321
$ws.Range("Q8").Value = 0.5276613635733334
$ws.Range("R8").Value = 4.748952272159999
$ws.Range("S8").Value = 0.0001699090139926458
$ws.Range("T8").Value = 0.0001699090139926458

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ptprc"
$ws.Range("C9").Value = "Mrc1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.04189333333333333
$ws.Range("H9").Value = 0.12568
$ws.Range("I9").Value = 0.0001717735044811748
$ws.Range("J9").Value = 0.0001717735044811748
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.07889166666666667
$ws.Range("N9").Value = 0.236675
$ws.Range("O9").Value = 0.006195566123830276
$ws.Range("P9").Value = 0.006195566123830277
$ws.Range("Q9").Value = 0.003305034888888889
$ws.Range("R9").Value = 0.02974531399999999
$ws.Range("S9").Value = 0.000001064234105335175
$ws.Range("T9").Value = 0.000001064234105335175

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ptprc"
$ws.Range("C10").Value = "Mrc1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.04189333333333333
$ws.Range("H10").Value = 0.12568
$ws.Range("I10").Value = 0.0001717735044811748
$ws.Range("J10").Value = 0.0001717735044811748
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.05932299999999999
$ws.Range("N10").Value = 0.177969
$ws.Range("O10").Value = 0.004658788243337701
$ws.Range("P10").Value = 0.004658788243337702
$ws.Range("Q10").Value = 0.002485238213333333
$ws.Range("R10").Value = 0.02236714392
$ws.Range("S10").Value = 0.000000800256383193813
$ws.Range("T10").Value = 0.0000008002563831938132
